$d = $word.ActiveDocument

# --- Paragraph 3 ("Adresse :") : merge split runs/proofErr into a single run ---
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertXML('<w:p w14:paraId="4875C844" w14:textId="203D648A" w:rsidR="00F53CF0" w:rsidRPr="001C7307" w:rsidRDefault="00551B13" w:rsidP="00CB3D9D"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Adresse :</w:t></w:r><w:r w:rsidR="00F53CF0" w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> 329 Rue Des Epoux Tramier, 84410, Bédoin</w:t></w:r></w:p>')

# --- Paragraph 4 ("Tel :") : drop the gramStart/gramEnd proofErr markers ---
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertXML('<w:p w14:paraId="75AE7A67" w14:textId="675B68CD" w:rsidR="009E0B57" w:rsidRPr="001C7307" w:rsidRDefault="00551B13" w:rsidP="00CB3D9D"><w:r w:rsidRPr="001C7307"><w:t>Tel :</w:t></w:r><w:r w:rsidR="00555429" w:rsidRPr="001C7307"><w:t xml:space="preserve"> 07 83 14 41 64</w:t></w:r></w:p>')

# --- Paragraph 5 ("Email :") : drop the gramStart/gramEnd proofErr markers ---
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertXML('<w:p w14:paraId="2BFA6955" w14:textId="01EAC39C" w:rsidR="00555429" w:rsidRPr="001C7307" w:rsidRDefault="00C55FED" w:rsidP="00CB3D9D"><w:pPr><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr></w:pPr><w:r w:rsidRPr="001C7307"><w:t>Email :</w:t></w:r><w:r w:rsidR="00555429" w:rsidRPr="001C7307"><w:t xml:space="preserve"> </w:t></w:r><w:hyperlink r:id="rId7" w:history="1"><w:r w:rsidR="00555429" w:rsidRPr="001C7307"><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>kenjiimbault@gmail.com</w:t></w:r></w:hyperlink></w:p>')

# --- Paragraph 6 ("Site web :") : drop the gramStart/gramEnd proofErr markers ---
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertXML('<w:p w14:paraId="671F4DEB" w14:textId="2AC8838A" w:rsidR="00D06B78" w:rsidRDefault="00C55FED" w:rsidP="00CB3D9D"><w:r w:rsidRPr="001C7307"><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="auto"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve">Site </w:t></w:r><w:r w:rsidR="0038710F" w:rsidRPr="001C7307"><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="auto"/><w:u w:val="none"/></w:rPr><w:t>web</w:t></w:r><w:r w:rsidR="009B76BE" w:rsidRPr="001C7307"><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="auto"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0038710F" w:rsidRPr="001C7307"><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="auto"/><w:u w:val="none"/></w:rPr><w:t>:</w:t></w:r><w:r w:rsidR="00755E6D" w:rsidRPr="001C7307"><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="auto"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:hyperlink r:id="rId8" w:history="1"><w:r w:rsidR="00D06B78" w:rsidRPr="009E1963"><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://kenjiimbault.tk/CV</w:t></w:r></w:hyperlink></w:p>')

# --- Paragraph 7 ("Permis de conduire") : drop the leading single-space Hyperlink run ---
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertXML('<w:p w14:paraId="092D4490" w14:textId="263B93C2" w:rsidR="00497F78" w:rsidRPr="00D63F13" w:rsidRDefault="00D06B78" w:rsidP="00CB3D9D"><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r w:rsidR="00CF506F" w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Permis de conduire</w:t></w:r><w:r w:rsidR="00497F78" w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00497F78" w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00497F78" w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00497F78" w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00497F78" w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00497F78" w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00497F78" w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00497F78" w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00CF506F" w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">Date </w:t></w:r><w:r w:rsidR="00C646DB" w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>d&#8217;obtention :</w:t></w:r><w:r w:rsidR="00497F78" w:rsidRPr="001C7307"><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:tab/><w:t>2020</w:t></w:r></w:p>')

# --- Paragraph with "ZCNP Foundation" : add "Zyxel " prefix run + drop a duplicate tab run ---
$zcnp = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "ZCNP Foundation*") {
        $zcnp = $d.Paragraphs.Item($i)
        break
    }
}
$zcnp.Range.InsertXML('<w:p w14:paraId="55700C98" w14:textId="453AA312" w:rsidR="001C7307" w:rsidRPr="00D63F13" w:rsidRDefault="002E6D4B" w:rsidP="001C7307"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Zyxel </w:t></w:r><w:r w:rsidRPr="006665CA"><w:t>ZCNP Foundation - v1.0 EN</w:t></w:r><w:r w:rsidRPr="006665CA"><w:tab/><w:t xml:space="preserve">ID: </w:t></w:r><w:r w:rsidRPr="00EA75D7"><w:t>6227d8406c53d13fbf3f4489</w:t></w:r><w:r w:rsidRPr="006665CA"><w:tab/></w:r><w:r><w:tab/></w:r><w:r w:rsidRPr="006665CA"><w:t>2022</w:t></w:r></w:p>')

# --- New paragraph: "Zyxel ZCNP Security - v1.0 EN" certification entry ---
$insertPoint = $d.Range($zcnp.Range.End, $zcnp.Range.End)
$insertPoint.InsertXML('<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Zyxel ZCNP Security - v1.0 EN</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t xml:space="preserve">ID: </w:t></w:r><w:r><w:t>62417ee3f34f3b3df92712bd</w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:t>2022</w:t></w:r></w:p>')
